# Batter vs team done
# Insert a new header row at the top of the sheet ("Team" / "Bowlers"),
# pushing the existing bowler/team rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row 1; existing rows 1..110 become 2..111.
$ws.Rows.Item(1).EntireRow.Insert()

# Populate the new header row.
$ws.Cells.Item(1, 1).Value = "Team"
$ws.Cells.Item(1, 2).Value = "Bowlers"

# Match the author's final selection/view state (cell B1 selected,
# scrolled back to the top of the sheet).
$ws.Range("B1").Select() | Out-Null
